$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# F6 previously held "Homepage3" (the lone testcase with that name);
# change it to "loginTest" (same group as F2) so the shared string
# "Homepage3" becomes unused and drops out of the table on save.
$ws.Range("F6").Value = "loginTest"

# Reflect the author's final selection in the sheet view.
$ws.Range("F6").Select()
